$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.54465389251709
$ws.Range("B1").Value = 4.123712062835693
$ws.Range("C1").Value = 3.566134214401245
$ws.Range("D1").Value = 4.34528636932373
$ws.Range("E1").Value = 4.800866603851318
